$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.115.52"
$ws.Range("E2").Value = "  +3.91%  "
$ws.Range("D3").Value = "2.496.16"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'322.63"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "'105.14"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D10").Value = "'37.64"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "'18.32"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").Value = "2.884.13"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "2.481.57"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "47.026.32"
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("D19").Value = "'12.62"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "'6.55"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").Value = "0.0₃0934"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'70.92"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").Value = "'251.27"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").Value = "'2.36"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").Value = "'26.18"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'10.17"
$ws.Range("E28").Value = "  +5.84%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").Value = "'35.39"
$ws.Range("E30").Value = "  +4.34%  "
$ws.Range("E31").Value = "  +4.63%  "
$ws.Range("D32").Value = "'49.58"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "'19.75"
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("D35").Value = "'0.0784"
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +3.02%  "

# Row 40/41: swap Stellar and Monero entries
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'122.53"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.111"
$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("D43").Value = "'21.63"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").Value = "'0.0294"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "1.953.37"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "'2.99"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "'1.78"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").Value = "'5.37"
$ws.Range("E50").Value = "  +12.96%  "
$ws.Range("D51").Value = "'78.73"
$ws.Range("E51").Value = "  +3.08%  "
